$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3315.1428
$ws.Range("I64").Value = 3150
$ws.Range("J64").Value = 3381.2
$ws.Range("K64").Value = 3150
$ws.Range("L64").Value = 3381.2
$ws.Range("M64").Value = -2902
$ws.Range("N64").Value = -3877.2
$ws.Range("H67").Value = 3315.1428
$ws.Range("I67").Value = 3150
$ws.Range("J67").Value = 3381.2
$ws.Range("K67").Value = 3150
$ws.Range("L67").Value = 3381.2
$ws.Range("M67").Value = -2292
$ws.Range("N67").Value = -5097.2
$ws.Range("H69").Value = 4291.1875
$ws.Range("I69").Value = 3650
$ws.Range("J69").Value = 4382.7856
$ws.Range("K69").Value = 10950
$ws.Range("L69").Value = 13148.3568
$ws.Range("M69").Value = -10076
$ws.Range("N69").Value = -14896.3568
$ws.Range("H72").Value = 4291.1875
$ws.Range("I72").Value = 3650
$ws.Range("J72").Value = 4382.7856
$ws.Range("K72").Value = 32850
$ws.Range("L72").Value = 39445.0704
$ws.Range("M72").Value = -28482
$ws.Range("N72").Value = -48181.0704
$ws.Range("H76").Value = 174634.61
$ws.Range("I76").Value = 249805.11
$ws.Range("J76").Value = 5501
$ws.Range("K76").Value = 249805.11
$ws.Range("L76").Value = 5501
$ws.Range("M76").Value = -249490.11
$ws.Range("N76").Value = -6131
$ws.Range("H79").Value = 174634.61
$ws.Range("I79").Value = 249805.11
$ws.Range("J79").Value = 5501
$ws.Range("K79").Value = 249805.11
$ws.Range("L79").Value = 5501
$ws.Range("M79").Value = -248713.11
$ws.Range("N79").Value = -7685
$ws.Range("H88").Value = 2010.5946
$ws.Range("I88").Value = 4383
$ws.Range("J88").Value = 1131.9259
$ws.Range("K88").Value = 4383
$ws.Range("L88").Value = 1131.9259
$ws.Range("M88").Value = -3977
$ws.Range("N88").Value = -1943.9259
$ws.Range("H91").Value = 2010.5946
$ws.Range("I91").Value = 4383
$ws.Range("J91").Value = 1131.9259
$ws.Range("K91").Value = 4383
$ws.Range("L91").Value = 1131.9259
$ws.Range("M91").Value = -2979
$ws.Range("N91").Value = -3939.9259
$ws.Range("H138").Value = 4255.32
$ws.Range("I138").Value = 1791.9375
$ws.Range("J138").Value = 4724.5356
$ws.Range("K138").Value = 5375.8125
$ws.Range("L138").Value = 14173.6068
$ws.Range("M138").Value = -235.8125
$ws.Range("N138").Value = -24453.6068

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 663.5417
$ws.Range("I2").Value = 584.1579
$ws.Range("J2").Value = 965.2
$ws.Range("K2").Value = 584.1579
$ws.Range("L2").Value = 965.2
$ws.Range("M2").Value = -471.1579
$ws.Range("N2").Value = -1191.2
$ws.Range("H116").Value = 663.5417
$ws.Range("I116").Value = 584.1579
$ws.Range("J116").Value = 965.2
$ws.Range("K116").Value = 584.1579
$ws.Range("L116").Value = 965.2
$ws.Range("M116").Value = 1709.8421
$ws.Range("N116").Value = -5553.2
$ws.Range("H132").Value = 26418.387
$ws.Range("I132").Value = 34094
$ws.Range("J132").Value = 3391.5454
$ws.Range("K132").Value = 102282
$ws.Range("L132").Value = 10174.6362
$ws.Range("M132").Value = -99752
$ws.Range("N132").Value = -15234.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 663.5417
$ws.Range("I3").Value = 584.1579
$ws.Range("J3").Value = 965.2
$ws.Range("K3").Value = 584.1579
$ws.Range("L3").Value = 965.2
$ws.Range("M3").Value = -470.1579
$ws.Range("N3").Value = -1193.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4578.3125
$ws.Range("I58").Value = 1023
$ws.Range("J58").Value = 12400
$ws.Range("K58").Value = 1023
$ws.Range("L58").Value = 12400
$ws.Range("M58").Value = -820
$ws.Range("N58").Value = -12806
$ws.Range("H132").Value = 2157.5833
$ws.Range("I132").Value = 1659.3077
$ws.Range("J132").Value = 2746.4546
$ws.Range("K132").Value = 4977.9231
$ws.Range("L132").Value = 8239.363799999999
$ws.Range("M132").Value = -2447.9231
$ws.Range("N132").Value = -13299.3638
$ws.Range("H136").Value = 4578.3125
$ws.Range("I136").Value = 1023
$ws.Range("J136").Value = 12400
$ws.Range("K136").Value = 3069
$ws.Range("L136").Value = 37200
$ws.Range("M136").Value = -519
$ws.Range("N136").Value = -42300

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 172.2
$ws.Range("I107").Value = 165.25
$ws.Range("J107").Value = 200
$ws.Range("K107").Value = 495.75
$ws.Range("L107").Value = 600
$ws.Range("M107").Value = 1424.25
$ws.Range("N107").Value = -4440

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2218.8462
$ws.Range("I126").Value = 1786.8
$ws.Range("K126").Value = 5360.4
$ws.Range("M126").Value = -2890.4
$ws.Range("H132").Value = 2920.9363
$ws.Range("I132").Value = 2996.973
$ws.Range("J132").Value = 2639.6
$ws.Range("K132").Value = 8990.919
$ws.Range("L132").Value = 7918.799999999999
$ws.Range("M132").Value = -6460.919
$ws.Range("N132").Value = -12978.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2564993.2
$ws.Range("I22").Value = 4167222.8
$ws.Range("J22").Value = 1426
$ws.Range("K22").Value = 4167222.8
$ws.Range("L22").Value = 1426
$ws.Range("M22").Value = -4166927.8
$ws.Range("N22").Value = -2016
$ws.Range("H27").Value = 2564993.2
$ws.Range("I27").Value = 4167222.8
$ws.Range("J27").Value = 1426
$ws.Range("K27").Value = 4167222.8
$ws.Range("L27").Value = 1426
$ws.Range("M27").Value = -4167115.8
$ws.Range("N27").Value = -1640
$ws.Range("H46").Value = 2042.8572
$ws.Range("I46").Value = 3800
$ws.Range("J46").Value = 725
$ws.Range("K46").Value = 3800
$ws.Range("L46").Value = 725
$ws.Range("M46").Value = -3612
$ws.Range("N46").Value = -1101
$ws.Range("H132").Value = 10210315
$ws.Range("I132").Value = 22738316
$ws.Range("J132").Value = 2313.1482
$ws.Range("K132").Value = 68214948
$ws.Range("L132").Value = 6939.444600000001
$ws.Range("M132").Value = -68212418
$ws.Range("N132").Value = -11999.4446
$ws.Range("H133").Value = 49500
$ws.Range("J133").Value = 49500
$ws.Range("L133").Value = 49500
$ws.Range("N133").Value = -54560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 41667056
$ws.Range("I113").Value = 425.25
$ws.Range("K113").Value = 1275.75
$ws.Range("M113").Value = 894.25
$ws.Range("H132").Value = 1592.6136
$ws.Range("I132").Value = 821.1429000000001
$ws.Range("J132").Value = 2942.6875
$ws.Range("K132").Value = 2463.4287
$ws.Range("L132").Value = 8828.0625
$ws.Range("M132").Value = 66.57129999999961
$ws.Range("N132").Value = -13888.0625
$ws.Range("H136").Value = 2970.0942
$ws.Range("I136").Value = 895.8823
$ws.Range("J136").Value = 6681.8423
$ws.Range("K136").Value = 2687.6469
$ws.Range("L136").Value = 20045.5269
$ws.Range("M136").Value = -137.6468999999997
$ws.Range("N136").Value = -25145.5269
